$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Restructure rows ---
# Table1 (Quadcopter parts) grows from 13 data rows to 16 data rows (+3):
#   new row inserted at row 3 ("Breadboard")
#   two new rows inserted at rows 8-9 ("3mm LED", "100R Resistor")
$ws.Range("A3:E3").EntireRow.Insert()
$ws.Range("A8:E9").EntireRow.Insert()

# Table2 (Receiver parts) shrinks from 7 data rows to 3 data rows (-4):
# after the inserts above, Table2's original data rows (20-26) now sit at rows 23-29:
#   23 Microprocessor (kept, becomes Arduino)
#   24 Resonator       (kept, becomes Breadboard)
#   25 3.3V Regulator  (removed)
#   26 10uF Capacitor  (removed)
#   27 10k Resistor    (removed)
#   28 2.4GHz Radio    (kept)
#   29 FTDI            (removed)
$ws.Range("A29:E29").EntireRow.Delete()
$ws.Range("A25:E27").EntireRow.Delete()

# --- Step 2: Resize the tables (ListObjects) to their new ranges ---
$lo1 = $ws.ListObjects.Item(1)
$lo2 = $ws.ListObjects.Item(2)
$lo1.Resize($ws.Range("A2:E19"))
$lo2.Resize($ws.Range("A22:E26"))

# --- Step 3: Ensure number formatting on the new rows' Cost cells ---
$ws.Range("C3").NumberFormat = $ws.Range("C4").NumberFormat
$ws.Range("C8:C9").NumberFormat = $ws.Range("C4").NumberFormat

# --- Step 4: Write final cell content (values / formulas) ---
# Row 1 (section header)
$ws.Range("A1").Value2 = "Quadcopter:"
# Row 2 (Table1 header)
$ws.Range("A2").Value2 = "Item"
$ws.Range("B2").Value2 = "Cat."
$ws.Range("C2").Value2 = "Cost"
$ws.Range("D2").Value2 = "From"
$ws.Range("E2").Value2 = "Comments"
# Row 3
$ws.Range("A3").Value2 = "Breadboard"
$ws.Range("B3").Value2 = "BREADBRD"
$ws.Range("C3").Value2 = 6.2
$ws.Range("D3").Value2 = "Futurlec"
$ws.Range("E3").Value2 = "Ok"
# Row 4
$ws.Range("A4").Value2 = "Microprocessor"
$ws.Range("B4").Value2 = "ATMEGA328P-PU"
$ws.Range("C4").Value2 = 2.31
$ws.Range("D4").Value2 = "Futurlec"
$ws.Range("E4").Value2 = "Ok"
# Row 5
$ws.Range("A5").Value2 = "Resonator"
$ws.Range("B5").Value2 = "RESON16M0P3"
$ws.Range("C5").Value2 = 0.21
$ws.Range("D5").Value2 = "Futurlec"
$ws.Range("E5").Value2 = "Ok"
# Row 6
$ws.Range("A6").Value2 = "10uF Capacitor"
$ws.Range("B6").Value2 = "C010U16E"
$ws.Range("C6").Value2 = 0.05
$ws.Range("D6").Value2 = "Futurlec"
$ws.Range("E6").Value2 = "Ok"
# Row 7
$ws.Range("A7").Value2 = "10k Resistor (x5)"
$ws.Range("B7").Value2 = "R010K14W"
$ws.Range("C7").Formula = "=0.11/10*5"
$ws.Range("D7").Value2 = "Futurlec"
$ws.Range("E7").Value2 = "Ok"
# Row 8 (new)
$ws.Range("A8").Value2 = "3mm LED"
$ws.Range("B8").Value2 = "LED3R"
$ws.Range("C8").Value2 = 0.08
$ws.Range("D8").Value2 = "Futurlec"
$ws.Range("E8").Value2 = "Ok"
# Row 9 (new)
$ws.Range("A9").Value2 = "100R Resistor"
$ws.Range("B9").Value2 = "R100R14W"
$ws.Range("C9").Formula = "=0.11/10"
$ws.Range("D9").Value2 = "Futurlec"
$ws.Range("E9").Value2 = "Ok"
# Row 10
$ws.Range("A10").Value2 = "Transistors (x4)"
$ws.Range("B10").Value2 = "BC337"
$ws.Range("C10").Formula = "=0.07*4"
$ws.Range("D10").Value2 = "Futurlec"
$ws.Range("E10").Value2 = "Ok"
# Row 11
$ws.Range("A11").Value2 = "Motors (x4)"
$ws.Range("B11").Value2 = "H107-A03"
$ws.Range("C11").Formula = "=11.82+0.99*11.82/10.57"
$ws.Range("D11").Value2 = "Futurlec"
$ws.Range("E11").Value2 = "Use Damo's motors to test"
# Row 12
$ws.Range("A12").Value2 = "2.4GHz Radio"
$ws.Range("B12").Value2 = "NRF24L01+"
$ws.Range("C12").Value2 = 1.44
$ws.Range("D12").Value2 = "eBay"
$ws.Range("E12").Value2 = "Ok"
# Row 13
$ws.Range("A13").Value2 = "3-axis acc + gyro"
$ws.Range("B13").Value2 = "MPU6050"
$ws.Range("C13").Value2 = 3.27
$ws.Range("D13").Value2 = "eBay"
$ws.Range("E13").Value2 = "Use 10DOF to test"
# Row 14
$ws.Range("A14").Value2 = "3-axis magnetometer"
$ws.Range("B14").Value2 = "HMC5883L"
$ws.Range("C14").Value2 = 2.48
$ws.Range("D14").Value2 = "eBay"
$ws.Range("E14").Value2 = "Use 10DOF to test"
# Row 15
$ws.Range("A15").Value2 = "Distance sensor"
$ws.Range("B15").Value2 = "HC-SR04"
$ws.Range("C15").Value2 = 1.79
$ws.Range("D15").Value2 = "eBay"
$ws.Range("E15").Value2 = "Ok"
# Row 16
$ws.Range("A16").Value2 = "Props"
$ws.Range("B16").Value2 = "H107-A02"
$ws.Range("C16").Value2 = 1.62
$ws.Range("D16").Value2 = "eBay"
$ws.Range("E16").Value2 = "Use HJ-998 props to test"
# Row 17
$ws.Range("A17").Value2 = "LiPo"
$ws.Range("B17").Value2 = "3.7V 350mAH 25C"
$ws.Range("C17").Formula = "=6.32/2"
$ws.Range("D17").Value2 = "eBay"
$ws.Range("E17").Value2 = "Ok"
# Row 18
$ws.Range("A18").Value2 = "LiPo Connector"
$ws.Range("B18").Value2 = "1 to 5 Charge Cable"
$ws.Range("C18").Formula = "=3.2/5"
$ws.Range("D18").Value2 = "eBay"
$ws.Range("E18").Value2 = "Ok"
# Row 19 (Table1 totals row)
$ws.Range("A19").Value2 = "Total"
$ws.Range("C19").Formula = "=SUBTOTAL(109,Table1[Cost])"

# Row 21 (section header)
$ws.Range("A21").Value2 = "Receiver:"
# Row 22 (Table2 header)
$ws.Range("A22").Value2 = "Item"
$ws.Range("B22").Value2 = "Cat."
$ws.Range("C22").Value2 = "Cost"
$ws.Range("D22").Value2 = "From"
$ws.Range("E22").Value2 = "Comments"
# Row 23
$ws.Range("A23").Value2 = "Arduino"
$ws.Range("B23").Value2 = "Arduino DK-UNO R3"
$ws.Range("C23").Value2 = 9.88
$ws.Range("D23").Value2 = "eBay"
$ws.Range("E23").Value2 = "Ok"
# Row 24
$ws.Range("A24").Value2 = "Breadboard"
$ws.Range("B24").Value2 = "BB390"
$ws.Range("C24").Value2 = 4.1
$ws.Range("D24").Value2 = "Futurlec"
$ws.Range("E24").Value2 = "Ok"
# Row 25
$ws.Range("A25").Value2 = "2.4GHz Radio"
$ws.Range("B25").Value2 = "NRF24L01+"
$ws.Range("C25").Value2 = 1.44
$ws.Range("D25").Value2 = "eBay"
$ws.Range("E25").Value2 = "Ok"
# Row 26 (Table2 totals row)
$ws.Range("A26").Value2 = "Total"
$ws.Range("C26").Formula = "=SUBTOTAL(109,Table2[Cost])"

# --- Step 5: fix up the sheet dimension / blank row 20 ---
$ws.Range("A20:E20").ClearContents()
